# #1031: Fix some typo in CodeList.rst
#
# The deck shows Spring XML snippets that reference
# "org.terasoluna.fw.common.codelist.*" CodeList classes. That package
# was renamed to "org.terasoluna.gfw.common.codelist.*" (the project
# moved from "terasoluna.fw" to "terasoluna.gfw"); slides 3 and 5 were
# already fixed, slides 1, 2 and 4 still had the stale "fw" package name.
# Also widens the CL_MONTH box on slide 2 so the now-longer class name
# line keeps fitting.

$p = $ppt.ActivePresentation

function Fix-CodeListClassName {
    param(
        [object]$Shape,
        [string]$OldClassName,
        [string]$NewClassName
    )

    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($OldClassName)
    if ($idx -lt 0) {
        return
    }

    $start = $idx + 1
    $len = $OldClassName.Length
    $chars = $tr.Characters($start, $len)
    $chars.Text = $NewClassName
}

# Slide 1: <bean id="CL_ORDERSTATUS" class="org.terasoluna.fw.common.codelist.SimpleMapCodeList">
$slide1 = $p.Slides.Item(1)
Fix-CodeListClassName $slide1.Shapes.Item(1) "org.terasoluna.fw.common.codelist.SimpleMapCodeList" "org.terasoluna.gfw.common.codelist.SimpleMapCodeList"

# Slide 2: <bean id="CL_MONTH"  class="org.terasoluna.fw.common.codelist.NumberRangeCodeList">
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(1)
Fix-CodeListClassName $shape2 "org.terasoluna.fw.common.codelist.NumberRangeCodeList" "org.terasoluna.gfw.common.codelist.NumberRangeCodeList"

# Widen the box slightly to account for the longer "gfw" class name text.
$shape2.Left = 46.5001
$shape2.Width = 385.4998

# Slide 4: <bean id="CL_AUTHORITIES" class="org.terasoluna.fw.common.codelist.JdbcCodeList">
$slide4 = $p.Slides.Item(4)
Fix-CodeListClassName $slide4.Shapes.Item(1) "org.terasoluna.fw.common.codelist.JdbcCodeList" "org.terasoluna.gfw.common.codelist.JdbcCodeList"
